$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3 (shifts old row 3 content down to row 4)
$ws.Rows.Item(3).Insert()

# The inserted row inherits extra formatted-but-empty cells (F3:H3) from the
# row above; remove them so only A3:D3 remain, matching the new data row.
$ws.Range("F3:H3").Clear()

# Fill in the new row 3 with the "recurring challenge" template entry.
# Cell-set order matters for shared-string table ordering.
$ws.Cells.Item(3, 1).Value = "Summary_Report"
$ws.Cells.Item(3, 3).Value = "One feature of the auto-generated template is that it pulls the most common recurring challenges across the agency. This text block summarizes how many times a recurring challenge has occurred for a given agency."
$ws.Cells.Item(3, 2).Value = "recurring_challenge_text"
$ws.Cells.Item(3, 4).Value = "**{challenge}** has been reported as challenge for the **{goal}** team in each of the last **{challenge count} quarters**."

# Set row heights to match the target layout
$ws.Rows.Item(3).RowHeight = 187.2
$ws.Rows.Item(4).RowHeight = 86.4

# Update the selected cell
$ws.Range("D3").Select()
